$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header "Save" in H1, matching the formatting of the other header
# cells (copy G1's format onto H1, same bold/centered/bordered header style)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# Add Save column values (0) for data rows 2-6
$ws.Range("H2:H6").Value = 0
